# Apply "Add data for 2022-05-01" changes:
# - Rename sheet from "Through 2022-04-22" to "Through 2022-04-23"
# - Update header label in I1 from "2022 (through 04-22)" to "2022 (through 04-23)"
# - Update April value I4 from 134 to 133
# - Update May value I5 from 92 to 103
# - Update Total value I14 from 528 to 538

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Through 2022-04-23"

$ws.Range("I1").Value = "2022 (through 04-23)"
$ws.Range("I4").Value = 133
$ws.Range("I5").Value = 103
$ws.Range("I14").Value = 538
